# Ajout reset esp32 et heartbeat
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# --- Update "Commentaires" column (E) on the test rows ---

# Row 10: Alimentation / short-circuit comment reworded with more detail about the fix
$ws.Range("E10").Value = @'
il y a un short-circuit entre le VUSB et le GND au niveau des via de support du connecteur usb-c. Problèeme réglé en perçant un trou dans le via dans les PCB.
'@.Trim()

# Row 12: Mécanique / poignée comment - typo fabricant fixed ("frabricant" -> "fabricant")
$ws.Range("E12").Value = @'
Le fabricant n'a pas mis de poignée comme désigné sur le PCB. La poignée est bien là sur Altium. Une plainte a été lancé.
'@.Trim()

# Rows 13, 14, 15, 16, 20: comments cleared out (removed from sharedStrings entirely)
$ws.Range("E13").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("E16").Value = ""

# Row 19: Communication i2c - new comment about esp32 boot order
$ws.Range("E19").Value = @'
Il faut que les esp32 boot après le PI pour activer le i2c des esp32
'@.Trim()

# Row 21: Alimentation intégrale - reworded phrasing
$ws.Range("E21").Value = @'
L'écran demande trop de courant si on l'alimente à l'aide du PI
'@.Trim()

# Row 22: Intégration - new expanded comment about pots / DEL / screen power
$ws.Range("E22").Value = @'
Les pots sont non fonctionnel. La DEL adressable au pin 10 est non fonctionnel. L'écran demande trop de courant si on l'alimente à l'aide du PI
'@.Trim()

# --- Row 35 of the "Modifications" table: new entry about SW resistor ordering ---
$ws.Range("C35").Value = @'
Les resistance des SW ne sont pas dans le bon ordre, ils sont dans l'ordre 2,1,4,3,6,5,8,7 au lieu de 1,2,3,4,5,6,7,8
'@.Trim()

# --- Restore scroll position / selection as left by the author ---
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
$ws.Range("C35").Select()
